# Fix a bug of null magic def useless.
# - Update the "Monster" sheet's L column (magic/tower id) for rows 4-50.
# - Make the "Monster" sheet the active tab (was "People"), and set its
#   active selection to L14.

$wb = $excel.ActiveWorkbook

# Sheet 2 = "怪物" (Monster) — holds the L-column id values that need fixing.
$ws = $wb.Worksheets.Item(2)

# New magic/tower-id values for L4:L50 (previously all 51018004).
$newValues = @{
    4  = 51018005
    5  = 51018005
    6  = 51018005
    7  = 51018006
    8  = 51018006
    9  = 51018006
    10 = 51018005
    11 = 51018006
    12 = 51018005
    13 = 51018005
    14 = 51018005
    15 = 51018005
    16 = 51018005
    17 = 51018005
    18 = 51018005
    19 = 51018007
    20 = 51018007
    21 = 51018005
    22 = 51018005
    23 = 51018007
    24 = 51018005
    25 = 51018006
    26 = 51018005
    27 = 51018005
    28 = 51018005
    29 = 51018005
    30 = 51018006
    31 = 51018006
    32 = 51018005
    33 = 51018005
    34 = 51018005
    35 = 51018005
    36 = 51018005
    37 = 51018005
    38 = 51018005
    39 = 51018007
    40 = 51018005
    41 = 51018005
    42 = 51018007
    43 = 51018007
    44 = 51018007
    45 = 51018005
    46 = 51018005
    47 = 51018005
    48 = 51018005
    49 = 51018005
    50 = 51018005
}

foreach ($row in $newValues.Keys) {
    $ws.Range("L$row").Value = $newValues[$row]
}

# Switch the active sheet from "People" to "Monster" and select L14, matching
# the saved view state (tabSelected / activeTab / selection) in the workbook.
$ws.Activate()
$ws.Range("L14").Select()
